$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.777.09'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.292.21'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '114.41'
$ws.Range("E5").Value = '  +19.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '268.50'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  +1.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.19'
$ws.Range("E10").Value = '  +5.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0938'
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.76'
$ws.Range("E12").Value = '  +12.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.61'
$ws.Range("E14").Value = '  +3.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.634.89'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.290.71'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.615.40'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("E19").Value = '  +2.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.53'
$ws.Range("E20").Value = '  +5.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.48'
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.48'
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.08'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.77'
$ws.Range("E24").Value = '  +6.40%  '
$ws.Range("E25").Value = '  +12.48%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.57'
$ws.Range("E27").Value = '  +4.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '42.13'
$ws.Range("E28").Value = '  +3.86%  '
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '176.60'
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0937'
$ws.Range("E32").Value = '  +5.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.59'
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.56'
$ws.Range("E34").Value = '  +3.70%  '
$ws.Range("E35").Value = '  +1.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.73'
$ws.Range("E36").Value = '  +9.25%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.85'
$ws.Range("E39").Value = '  +12.32%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.43'
$ws.Range("E40").Value = '  +5.29%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.85'
$ws.Range("E41").Value = '  +13.17%  '
$ws.Range("E42").Value = '  +3.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.89'
$ws.Range("E43").Value = '  +11.86%  '
$ws.Range("E44").Value = '  +7.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.06'
$ws.Range("E45").Value = '  +16.19%  '
$ws.Range("E46").Value = '  -0.19%  '
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.68'
$ws.Range("E48").Value = '  +5.77%  '
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.23'
$ws.Range("E50").Value = '  +3.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.449'
$ws.Range("E51").Value = '  +3.86%  '
